$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "date-looking" text value without letting Excel's
# autodetect convert it into a real date serial number + date format.
function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2: work/3000/2025-08-01 -> friend/343/2025-08-01
$ws.Range("A2").Value = "friend"
$ws.Range("B2").Value = 343
Set-TextValue "C2" "2025-08-01"

# Row 3: Prize/2000/2025-08-01 -> friend/500/2025-07-30
$ws.Range("A3").Value = "friend"
$ws.Range("B3").Value = 500
Set-TextValue "C3" "2025-07-30"

# Row 4: Hackathon/3000/2025-07-29 -> Interest/300/2025-07-29
$ws.Range("A4").Value = "Interest"
$ws.Range("B4").Value = 300
Set-TextValue "C4" "2025-07-29"

# Row 5: Friend/500/2025-07-28 -> Pocket money/5000/2025-07-27
$ws.Range("A5").Value = "Pocket money"
$ws.Range("B5").Value = 5000
Set-TextValue "C5" "2025-07-27"

# Row 6: Pocket money/5000/2025-07-27 -> Gift/500/2025-07-20
$ws.Range("A6").Value = "Gift"
$ws.Range("B6").Value = 500
Set-TextValue "C6" "2025-07-20"

# Row 7: Interest/200/2025-07-23 -> Salary/1000/2025-06-18
$ws.Range("A7").Value = "Salary"
$ws.Range("B7").Value = 1000
Set-TextValue "C7" "2025-06-18"

# Row 8: Gift/500/2025-07-20 -> removed entirely (row 8 deleted, shrinks
# the used range/dimension down to A1:C7)
$ws.Rows("8:8").Delete()
